$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive a new checkmark value (shared string "✓", same text
# used throughout the sheet for completed weeks).
$checkCells = @("U12", "V13", "U14", "V15", "U16", "W16", "X16")

# Cells whose style moves to the "checked" look (centered, Times New Roman,
# color FF0A1829) but keep no value (still blank).
$styleOnlyCells = @("V12", "U13")

function Set-CheckedStyle($range) {
    $range.Font.Name = "Times New Roman"
    $range.Font.Size = 12
    $range.Font.Color = 2693130
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4108
}

foreach ($addr in $checkCells) {
    $r = $ws.Range($addr)
    Set-CheckedStyle $r
    $r.Value = "✓"
}

foreach ($addr in $styleOnlyCells) {
    $r = $ws.Range($addr)
    Set-CheckedStyle $r
}

# Sheet view's active selection moved from M5 to T6.
$ws.Range("T6").Select()
